$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values while preserving their original "text" cell type.
# Excel auto-coerces numeric-looking strings (e.g. "0.999", "7.02") into
# real numbers on plain .Value assignment, which would both change the
# stored cell type and introduce floating point artifacts. Forcing a Text
# number format before the write keeps the literal string, and resetting
# the style back to Normal afterwards avoids leaving a stray style index
# on the cell (matching the un-styled inline strings in the source file).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '58.853.83'
Set-TextValue $ws.Range('E2') '  -2.18%  '
Set-TextValue $ws.Range('D3') '2.652.59'
Set-TextValue $ws.Range('E3') '  -0.62%  '
Set-TextValue $ws.Range('E4') '  -0.07%  '
Set-TextValue $ws.Range('D5') '522.93'
Set-TextValue $ws.Range('E5') '  +0.55%  '
Set-TextValue $ws.Range('D6') '144.28'
Set-TextValue $ws.Range('E6') '  -0.88%  '
Set-TextValue $ws.Range('D7') '0.999'
Set-TextValue $ws.Range('E7') '  +0.17%  '
Set-TextValue $ws.Range('D8') '0.572'
Set-TextValue $ws.Range('E8') '  -1.01%  '
Set-TextValue $ws.Range('D9') '7.02'
Set-TextValue $ws.Range('E9') '  +9.20%  '
Set-TextValue $ws.Range('E10') '  -2.38%  '
Set-TextValue $ws.Range('D11') '0.334'
Set-TextValue $ws.Range('E11') '  -1.70%  '
Set-TextValue $ws.Range('E12') '  +1.45%  '
Set-TextValue $ws.Range('D13') '3.118.08'
Set-TextValue $ws.Range('E13') '  -0.56%  '
Set-TextValue $ws.Range('D14') '58.865.70'
Set-TextValue $ws.Range('E14') '  -2.18%  '
Set-TextValue $ws.Range('D15') '21.04'
Set-TextValue $ws.Range('E15') '  -0.88%  '
Set-TextValue $ws.Range('B16') 'WrappedEther'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D16') '2.678.56'
Set-TextValue $ws.Range('E16') '  -2.89%  '
Set-TextValue $ws.Range('B17') 'ShibaInu'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D17') '0.0000136'
Set-TextValue $ws.Range('E17') '  -1.62%  '
Set-TextValue $ws.Range('D18') '338.62'
Set-TextValue $ws.Range('E18') '  -3.04%  '
Set-TextValue $ws.Range('D19') '4.36'
Set-TextValue $ws.Range('E19') '  -3.94%  '
Set-TextValue $ws.Range('D20') '10.36'
Set-TextValue $ws.Range('E20') '  -1.44%  '
Set-TextValue $ws.Range('D21') '6.35'
Set-TextValue $ws.Range('E21') '  +0.90%  '
Set-TextValue $ws.Range('E22') '  +0.12%  '
Set-TextValue $ws.Range('D23') '63.78'
Set-TextValue $ws.Range('E23') '  +1.78%  '
Set-TextValue $ws.Range('E24') '  -0.73%  '
Set-TextValue $ws.Range('D25') '0.166'
Set-TextValue $ws.Range('E25') '  -1.01%  '
Set-TextValue $ws.Range('E26') '  +0.40%  '
Set-TextValue $ws.Range('D27') '0.0₃0800'
Set-TextValue $ws.Range('E27') '  -1.17%  '
Set-TextValue $ws.Range('D28') '7.06'
Set-TextValue $ws.Range('E28') '  -2.32%  '
Set-TextValue $ws.Range('E29') '  -2.30%  '
Set-TextValue $ws.Range('E30') '  -0.04%  '
Set-TextValue $ws.Range('E31') '  +0.19%  '
Set-TextValue $ws.Range('D32') '18.83'
Set-TextValue $ws.Range('E32') '  -0.92%  '
Set-TextValue $ws.Range('D33') '149.43'
Set-TextValue $ws.Range('E33') '  +0.71%  '
Set-TextValue $ws.Range('E34') '  -3.51%  '
Set-TextValue $ws.Range('E35') '  -2.77%  '
Set-TextValue $ws.Range('D36') '0.890'
Set-TextValue $ws.Range('E36') '  -6.07%  '
Set-TextValue $ws.Range('D37') '0.868'
Set-TextValue $ws.Range('E37') '  -0.17%  '
Set-TextValue $ws.Range('D38') '36.73'
Set-TextValue $ws.Range('E38') '  +0.18%  '
Set-TextValue $ws.Range('E39') '  -5.87%  '
Set-TextValue $ws.Range('E40') '  -2.69%  '
Set-TextValue $ws.Range('D41') '0.616'
Set-TextValue $ws.Range('E41') '  +1.53%  '
Set-TextValue $ws.Range('E42') '  +0.34%  '
Set-TextValue $ws.Range('E43') '  +0.10%  '
Set-TextValue $ws.Range('D44') '275.38'
Set-TextValue $ws.Range('E44') '  -2.03%  '
Set-TextValue $ws.Range('D45') '0.0969'
Set-TextValue $ws.Range('E45') '  -1.98%  '
Set-TextValue $ws.Range('E46') '  +1.98%  '
Set-TextValue $ws.Range('E47') '  -1.28%  '
Set-TextValue $ws.Range('D48') '2.039.14'
Set-TextValue $ws.Range('E48') '  -3.83%  '
Set-TextValue $ws.Range('E49') '  -2.76%  '
Set-TextValue $ws.Range('B50') 'InjectiveProtocol'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D50') '18.89'
Set-TextValue $ws.Range('E50') '  -0.50%  '
Set-TextValue $ws.Range('B51') 'VeChain'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D51') '0.0228'
Set-TextValue $ws.Range('E51') '  -2.74%  '
